# Weekly data refresh: insert the two newest price records at the top of
# the data block (rows 109-110), pushing the existing records down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 109, shifting the
# existing rows 109-130 down to 111-132.
$ws.Rows("109:110").Insert()

# --- New row 109 ---
$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44463
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100104
$ws.Range("H109").Value = "Frutos de pepita"
$ws.Range("I109").Value = 100104005
$ws.Range("J109").Value = "Pera"
$ws.Range("K109").Value = "Packham's Triumph"
$ws.Range("L109").Value = "Especial"
$ws.Range("M109").Value = 30
$ws.Range("N109").Value = 11000
$ws.Range("O109").Value = 11000
$ws.Range("P109").Value = 11000
$ws.Range("Q109").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R109").Value = "Provincia de Curicó"
$ws.Range("S109").Value = 688
$ws.Range("T109").Value = 16

# --- New row 110 ---
$ws.Range("A110").Value = 7
$ws.Range("B110").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C110").Value = "Ñuble"
$ws.Range("D110").Value = 44463
$ws.Range("E110").Value = 16
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100104
$ws.Range("H110").Value = "Frutos de pepita"
$ws.Range("I110").Value = 100104005
$ws.Range("J110").Value = "Pera"
$ws.Range("K110").Value = "Packham's Triumph"
$ws.Range("L110").Value = "Primera"
$ws.Range("M110").Value = 60
$ws.Range("N110").Value = 9000
$ws.Range("O110").Value = 10000
$ws.Range("P110").Value = 9500
$ws.Range("Q110").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R110").Value = "Provincia de Curicó"
$ws.Range("S110").Value = 594
$ws.Range("T110").Value = 16
